$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing header cell
# onto the three new header cells before setting their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins / Losses / Ties) for every player row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 61
    $ws.Cells.Item($r, 31).Value = 101
    $ws.Cells.Item($r, 32).Value = 0
}
